$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 127.3992563333333
$ws.Range("N2").Value = 382.197769
$ws.Range("O2").Value = 0.4838549810199306
$ws.Range("P2").Value = 0.4838549810199307
$ws.Range("Q2").Value = 93.2751531923561
$ws.Range("R2").Value = 839.4763787312049
$ws.Range("S2").Value = 0.02474202333198028
$ws.Range("T2").Value = 0.02474202333198029
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("M3").Value = 59.36586533333332
$ws.Range("O3").Value = 0.2254681108101269
$ws.Range("P3").Value = 0.2254681108101269
$ws.Range("Q3").Value = 43.4646193606911
$ws.Range("R3").Value = 391.1815742462199
$ws.Range("S3").Value = 0.01152935792150476
$ws.Range("T3").Value = 0.01152935792150477
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 16.63275166666667
$ws.Range("N4").Value = 49.898255
$ws.Range("O4").Value = 0.06317022542837675
$ws.Range("P4").Value = 0.06317022542837675
$ws.Range("Q4").Value = 12.17764141149722
$ws.Range("R4").Value = 109.598772703475
$ws.Range("S4").Value = 0.003230222386345489
$ws.Range("T4").Value = 0.003230222386345489
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("M5").Value = 59.90262233333334
$ws.Range("N5").Value = 179.707867
$ws.Range("O5").Value = 0.2275066827415657
$ws.Range("P5").Value = 0.2275066827415658
$ws.Range("Q5").Value = 43.85760510364612
$ws.Range("R5").Value = 394.718445932815
$ws.Range("S5").Value = 0.01163360071380848
$ws.Range("T5").Value = 0.01163360071380849
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 127.3992563333333
$ws.Range("N6").Value = 382.197769
$ws.Range("O6").Value = 0.4838549810199306
$ws.Range("P6").Value = 0.4838549810199307
$ws.Range("Q6").Value = 1307.095901166149
$ws.Range("R6").Value = 11763.86311049534
$ws.Range("S6").Value = 0.3467182435722757
$ws.Range("T6").Value = 0.3467182435722758
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("M7").Value = 59.36586533333332
$ws.Range("O7").Value = 0.2254681108101269
$ws.Range("P7").Value = 0.2254681108101269
$ws.Range("Q7").Value = 609.084240204303
$ws.Range("S7").Value = 0.161564746521492
$ws.Range("T7").Value = 0.161564746521492
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 16.63275166666667
$ws.Range("N8").Value = 49.898255
$ws.Range("O8").Value = 0.06317022542837675
$ws.Range("P8").Value = 0.06317022542837675
$ws.Range("Q8").Value = 170.6493597712322
$ws.Range("R8").Value = 1535.84423794109
$ws.Range("S8").Value = 0.04526618608001745
$ws.Range("T8").Value = 0.04526618608001745
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("M9").Value = 59.90262233333334
$ws.Range("N9").Value = 179.707867
$ws.Range("O9").Value = 0.2275066827415657
$ws.Range("P9").Value = 0.2275066827415658
$ws.Range("Q9").Value = 614.5912807853452
$ws.Range("R9").Value = 5531.321527068107
$ws.Range("S9").Value = 0.163025535615725
$ws.Range("T9").Value = 0.163025535615725
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 127.3992563333333
$ws.Range("N10").Value = 382.197769
$ws.Range("O10").Value = 0.4838549810199306
$ws.Range("P10").Value = 0.4838549810199307
$ws.Range("Q10").Value = 327.1613327304188
$ws.Range("R10").Value = 2944.451994573769
$ws.Range("S10").Value = 0.08678231073011136
$ws.Range("T10").Value = 0.08678231073011138
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("M11").Value = 59.36586533333332
$ws.Range("O11").Value = 0.2254681108101269
$ws.Range("P11").Value = 0.2254681108101269
$ws.Range("Q11").Value = 152.4515619646218
$ws.Range("R11").Value = 1372.064057681596
$ws.Range("S11").Value = 0.04043906628967747
$ws.Range("T11").Value = 0.04043906628967748
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 16.63275166666667
$ws.Range("N12").Value = 49.898255
$ws.Range("O12").Value = 0.06317022542837675
$ws.Range("P12").Value = 0.06317022542837675
$ws.Range("Q12").Value = 42.71291182425056
$ws.Range("R12").Value = 384.416206418255
$ws.Range("S12").Value = 0.01132996113930831
$ws.Range("T12").Value = 0.01132996113930831
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("M13").Value = 59.90262233333334
$ws.Range("N13").Value = 179.707867
$ws.Range("O13").Value = 0.2275066827415657
$ws.Range("P13").Value = 0.2275066827415658
$ws.Range("Q13").Value = 153.8299541195408
$ws.Range("R13").Value = 1384.469587075867
$ws.Range("S13").Value = 0.04080469646760164
$ws.Range("T13").Value = 0.04080469646760165
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 127.3992563333333
$ws.Range("N14").Value = 382.197769
$ws.Range("O14").Value = 0.4838549810199306
$ws.Range("P14").Value = 0.4838549810199307
$ws.Range("Q14").Value = 96.55640597205866
$ws.Range("R14").Value = 869.0076537485279
$ws.Range("S14").Value = 0.02561240338556328
$ws.Range("T14").Value = 0.02561240338556329
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("M15").Value = 59.36586533333332
$ws.Range("O15").Value = 0.2254681108101269
$ws.Range("P15").Value = 0.2254681108101269
$ws.Range("Q15").Value = 44.99362679959466
$ws.Range("R15").Value = 404.9426411963519
$ws.Range("S15").Value = 0.01193494007745263
$ws.Range("T15").Value = 0.01193494007745263
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 16.63275166666667
$ws.Range("N16").Value = 49.898255
$ws.Range("O16").Value = 0.06317022542837675
$ws.Range("P16").Value = 0.06317022542837675
$ws.Range("Q16").Value = 12.60602901917333
$ws.Range("R16").Value = 113.45426117256
$ws.Range("S16").Value = 0.003343855822705496
$ws.Range("T16").Value = 0.003343855822705496
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("M17").Value = 59.90262233333334
$ws.Range("N17").Value = 179.707867
$ws.Range("O17").Value = 0.2275066827415657
$ws.Range("P17").Value = 0.2275066827415658
$ws.Range("Q17").Value = 45.40043707692266
$ws.Range("R17").Value = 408.603933692304
$ws.Range("S17").Value = 0.01204284994443062
$ws.Range("T17").Value = 0.01204284994443062
